$wb = $excel.ActiveWorkbook

# "Generate Report for Handback" - refresh the handoff/handback timestamps
# for the zh-cn and de-de report rows.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-18 10:49:16"
$wsZhCn.Range("H2").Value = "2016-03-18 10:49:34"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-18 10:49:19"
$wsDeDe.Range("H2").Value = "2016-03-18 10:49:39"
